$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.590.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -7.37%  '

$ws.Range("D3").Value = "'1.685.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.57%  '

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'217.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.34%  '

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D6").Value = "'0.5007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -15.73%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = "'0.2611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.34%  '

$ws.Range("D9").Value = "'21.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.94%  '

$ws.Range("D10").Value = "'0.06205"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.22%  '

$ws.Range("D11").Value = "'0.07286"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.43%  '

$ws.Range("D12").Value = "'1.664.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.34%  '

$ws.Range("D13").Value = "'4.446"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.47%  '

$ws.Range("D14").Value = "'0.5764"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.54%  '

$ws.Range("D15").Value = "'1.911.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.71%  '

$ws.Range("D16").Value = "'0.000008218"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -11.83%  '

$ws.Range("D17").Value = "'64.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -14.55%  '

$ws.Range("D18").Value = "'26.573.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.28%  '

$ws.Range("D19").Value = "'5.002"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.00%  '

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").Value = "'10.78"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'185.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.88%  '

$ws.Range("D23").Value = "'6.205"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -9.68%  '

$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").Value = "'144.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.24%  '

$ws.Range("D26").Value = "'7.498"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.81%  '

$ws.Range("D27").Value = "'0.1135"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.84%  '

$ws.Range("D28").Value = "'15.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.62%  '

$ws.Range("E29").Value = '  -8.64%  '

$ws.Range("D30").Value = "'0.05730"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.27%  '

$ws.Range("E31").Value = '  -7.27%  '

$ws.Range("D32").Value = "'3.483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.04%  '

$ws.Range("D33").Value = "'3.479"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.37%  '

$ws.Range("D34").Value = "'1.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.67%  '

$ws.Range("D35").Value = "'1.008"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.48%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'0.5922"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.70%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = "'2.366"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.07%  '

$ws.Range("E38").Value = '  -3.08%  '

$ws.Range("D39").Value = "'0.01590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.29%  '

$ws.Range("D40").Value = "'1.073.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.37%  '

$ws.Range("D41").Value = "'5.885"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -9.11%  '

$ws.Range("D42").Value = "'0.8542"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.79%  '

$ws.Range("D43").Value = "'1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.74%  '

$ws.Range("D44").Value = "'98.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.31%  '

$ws.Range("D45").Value = "'1.839.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.48%  '

$ws.Range("D46").Value = "'56.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.92%  '

$ws.Range("D47").Value = "'0.00000000106"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.62%  '

$ws.Range("E48").Value = '  -0.49%  '

$ws.Range("D49").Value = "'8.036"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.64%  '

$ws.Range("E50").Value = '  -3.91%  '

$ws.Range("D51").Value = "'0.05197"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.09%  '
